# The workbook's "Feuil1" sheet is a data-dictionary / record layout table.
# Row 19 ("N° séquentiel de séjour") had its variable name ("nom", column F)
# mis-labelled as "NOSEJHAD" (a name that belongs to the HAD format).
# This commit ("dans les psy ano aussi") fixes it to the correct name
# "NOSEQSEJ", mirroring the same fix already done in the other ANO formats.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Fix the mislabeled variable name in column F (nom) for row 19.
$ws.Range("F19").Value = "NOSEQSEJ"

# Reflect the cell that was left selected/active after the edit.
$ws.Range("F20").Select() | Out-Null

# Column A was widened (best-fit) to accommodate the long "libelle" text.
# (87.5 characters wide, as stored in the saved XML's <col width="...">.)
$ws.Columns.Item(1).ColumnWidth = 86.66666666666667
